$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (OOXML <col> "width" character units) ---
# Excel's COM ColumnWidth is the "characters" width; the serialized <col width>
# attribute adds ~5/6 of a character of internal cell padding on top of that,
# so we dial the COM value back by 5/6 to land exactly on the target width.
$pad = 5 / 6
$ws.Columns.Item(3).ColumnWidth = 76 - $pad
$ws.Columns.Item(4).ColumnWidth = 54 - $pad
$ws.Columns.Item(6).ColumnWidth = 17 - $pad
$ws.Columns.Item(7).ColumnWidth = 16 - $pad
$ws.Columns.Item(8).ColumnWidth = 34 - $pad

# --- Row data for rows 2..16 (A..H) ---
$data = @(
    @("1328652", "https://aiesec.org/opportunity/global-talent/1328652", "AI Engineer", "Novi Sad, Serbia", "No", "1 applicant", "9 - 12 Weeks", "DataDrill"),
    @("1328624", "https://aiesec.org/opportunity/global-talent/1328624", "Sales Representative", "İstanbul, Türkiye", "No", "7 applicants", "6 - 18 Months", "Saphiredent"),
    @("1328613", "https://aiesec.org/opportunity/global-talent/1328613", "Sales Manager Intern", "Adana, Reşatbey, Seyhan/Adana, Türkiye", "No", "8 applicants", "9 - 12 Weeks", "HYFROST"),
    @("1328604", "https://aiesec.org/opportunity/global-talent/1328604", "Market Analysis and Improvement Proposals Based on Call Center Operations", "日本、東京都東京", "No", "2 applicants", "9 - 12 Weeks", "CENTRIC Co., Ltd."),
    @("1328537", "https://aiesec.org/opportunity/global-talent/1328537", "French Preschool Coordinator", "Cairo, Cairo Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Talent School - October Gardens"),
    @("1328536", "https://aiesec.org/opportunity/global-talent/1328536", "English Preschool Coordinator", "Cairo, Cairo Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Talent School - October Gardens"),
    @("1328533", "https://aiesec.org/opportunity/global-talent/1328533", "WordPress Developer", "Cairo, Cairo Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Think Advertising"),
    @("1328518", "https://aiesec.org/opportunity/global-talent/1328518", "Social Media Specialist", "Cairo, Cairo Governorate, Egypt", "No", "0 applicants", "6 - 18 Months", "Tech-Pix"),
    @("1328517", "https://aiesec.org/opportunity/global-talent/1328517", "Videographer", "Cairo, Cairo Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Tech-Pix"),
    @("1327541", "https://aiesec.org/opportunity/global-talent/1327541", "Software Developer", "União das freguesias de Cascais e Estoril, Portugal", "No", "92 applicants", "3 - 6 Months", "Dark Cloud"),
    @("1327539", "https://aiesec.org/opportunity/global-talent/1327539", "Graphic Designer", "União das freguesias de Cascais e Estoril, Portugal", "No", "34 applicants", "3 - 6 Months", "Dark Cloud"),
    @("1327156", "https://aiesec.org/opportunity/global-talent/1327156", "Security Master Domain Intern", "Heerlen, Netherlands", "No", "64 applicants", "6 - 18 Months", "APG Heerlen"),
    @("1325702", "https://aiesec.org/opportunity/global-talent/1325702", "Guest Relations Executive and Waitress", "Colombo, Sri Lanka", "No", "12 applicants", "3 - 6 Months", "Indian Kitchen PVT LTD"),
    @("1325700", "https://aiesec.org/opportunity/global-talent/1325700", "Business Analyst and Executive Secretary", "Colombo, Sri Lanka", "No", "26 applicants", "3 - 6 Months", "Indian Kitchen PVT LTD"),
    @("1316788", "https://aiesec.org/opportunity/global-talent/1316788", "Travel Coordinator", "Mexico City, CDMX, Mexico", "No", "123 applicants", "6 - 18 Months", "Ikan Experience"),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]

    # Column A holds numeric-looking opportunity IDs that must remain TEXT
    # (as in the source data), so prefix with an apostrophe text qualifier
    # to stop Excel from auto-converting them to numbers.
    $ws.Cells.Item($row, 1).Value = "'" + $values[0]

    for ($j = 1; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}
